$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.003783345222473
$ws.Range("B1").Value = 2.113560438156128
$ws.Range("C1").Value = 6.720728397369385
$ws.Range("D1").Value = 1.894567608833313
$ws.Range("E1").Value = 1.370669007301331
